# Fix missing prefab pipeline:
#  - "RecontraBird" was a placeholder name; rename it to "Bird" (the same
#    display name already used by the "Bird view" prefab row), which also
#    collapses the now-duplicate standalone "Bird" shared string.
#  - Row 5 on the Animals sheet was an accidental duplicate of row 4
#    (same IDS/Life/Food/Name); clear its contents, leaving an empty row.
#  - Row 4's ID cell (A4) picks up the alignment formatting already used
#    elsewhere in the sheet (e.g. B2), matching style "2" (readingOrder=0).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Animals")

# Rename the animal at row 4 from "RecontraBird" to "Bird".
$ws.Range("A4").Value2 = "Bird"

# Pick up the alignment format already used on the sheet (B2 carries the
# "readingOrder=0" alignment applied via Format Cells), so A4 matches it.
$ws.Range("B2").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 5 duplicated row 4's data; clear it out, keeping the row/style shell.
$ws.Range("A5:D5").ClearContents()
